# "Generate Report for Handback"
#
# The localization-status report is regenerated: the status for the
# zh-cn and de-de handback rows flips from "Ready for handoff" to
# "Handed back: in sync with en-US", the Latest Handback DateTime
# stamps are refreshed, and the (now stale/resolved) handback-version
# error message is cleared out. Column widths on the Status / Error
# Detail columns are refreshed to fit the new text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet ------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Widen the zh-cn / de-de status columns to fit the longer text.
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---- zh-cn sheet -----------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-28 08:48:12"
$zhcn.Range("P2").Value = ""

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334

# ---- de-de sheet -----------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-28 08:48:19"
$dede.Range("P2").Value = ""

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333334
